$wb = $excel.ActiveWorkbook

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 868.7143
$ws.Range("I28").Value = 1073
$ws.Range("J28").Value = 358
$ws.Range("K28").Value = 1073
$ws.Range("L28").Value = 358
$ws.Range("M28").Value = -588
$ws.Range("N28").Value = -1328

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 190
$ws.Range("I33").Value = 190
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 190
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 39
$ws.Range("N33").Value = $null

# ALC row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 2850.125
$ws.Range("I41").Value = 3300.1667
$ws.Range("J41").Value = 1500
$ws.Range("K41").Value = 3300.1667
$ws.Range("L41").Value = 1500
$ws.Range("M41").Value = -2860.1667
$ws.Range("N41").Value = -2380

# ALC row 69
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 5000
$ws.Range("I69").Value = 5000
$ws.Range("J69").Value = 5000
$ws.Range("K69").Value = 15000
$ws.Range("L69").Value = 15000
$ws.Range("M69").Value = -14126
$ws.Range("N69").Value = -16748

# ALC row 72
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 5000
$ws.Range("I72").Value = 5000
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 45000
$ws.Range("L72").Value = 45000
$ws.Range("M72").Value = -40632
$ws.Range("N72").Value = -53736

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4351.3335
$ws.Range("J86").Value = 4351.3335
$ws.Range("L86").Value = 4351.3335
$ws.Range("N86").Value = -6597.3335

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4351.3335
$ws.Range("J89").Value = 4351.3335
$ws.Range("L89").Value = 21756.6675
$ws.Range("N89").Value = -32988.6675

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").Value = $null

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3997
$ws.Range("I106").Value = 3997
$ws.Range("K106").Value = 3997
$ws.Range("M106").Value = -3366

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1380.6666
$ws.Range("I107").Value = 806.9
$ws.Range("K107").Value = 806.9
$ws.Range("M107").Value = 1113.1

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7597.8335
$ws.Range("I132").Value = 6646.75
$ws.Range("K132").Value = 19940.25
$ws.Range("M132").Value = -17410.25

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4946
$ws.Range("I2").Value = 4946
$ws.Range("K2").Value = 4946
$ws.Range("M2").Value = -4833

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1075
$ws.Range("I110").Value = 650
$ws.Range("K110").Value = 650
$ws.Range("M110").Value = 1395

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 4946
$ws.Range("I116").Value = 4946
$ws.Range("K116").Value = 4946
$ws.Range("M116").Value = -2652

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4946
$ws.Range("I3").Value = 4946
$ws.Range("K3").Value = 4946
$ws.Range("M3").Value = -4832

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4105.364
$ws.Range("I99").Value = 4105.364
$ws.Range("K99").Value = 4105.364
$ws.Range("M99").Value = -2607.364

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4222.2
$ws.Range("I107").Value = 4222.2
$ws.Range("K107").Value = 4222.2
$ws.Range("M107").Value = -2302.2

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5493.5
$ws.Range("I134").Value = 990
$ws.Range("K134").Value = 2970
$ws.Range("M134").Value = -435

# CRP row 54
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 83
$ws.Range("I54").Value = 83
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 83
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = 575
$ws.Range("N54").Value = $null

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6166.6665
$ws.Range("I70").Value = 5500
$ws.Range("J70").Value = 7500
$ws.Range("K70").Value = 5500
$ws.Range("L70").Value = 7500
$ws.Range("M70").Value = -5230
$ws.Range("N70").Value = -8040

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6166.6665
$ws.Range("I73").Value = 5500
$ws.Range("J73").Value = 7500
$ws.Range("K73").Value = 5500
$ws.Range("L73").Value = 7500
$ws.Range("M73").Value = -4564
$ws.Range("N73").Value = -9372

# GSM row 93
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = $null

# GSM row 98
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 46749.25
$ws.Range("J98").Value = 46749.25
$ws.Range("L98").Value = 46749.25
$ws.Range("N98").Value = -52739.25

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5250
$ws.Range("I7").Value = 5250
$ws.Range("K7").Value = 5250
$ws.Range("M7").Value = -5138

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6156.8
$ws.Range("I40").Value = 6156.8
$ws.Range("K40").Value = 6156.8
$ws.Range("M40").Value = -6020.8

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4928.4287
$ws.Range("I122").Value = 4999.75
$ws.Range("J122").Value = 4833.3335
$ws.Range("K122").Value = 14999.25
$ws.Range("L122").Value = 14500.0005
$ws.Range("M122").Value = -12549.25
$ws.Range("N122").Value = -19400.0005

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5250
$ws.Range("I126").Value = 5250
$ws.Range("K126").Value = 15750
$ws.Range("M126").Value = -13280

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4668.615
$ws.Range("I132").Value = 4619.7
$ws.Range("J132").Value = 4831.6665
$ws.Range("K132").Value = 13859.1
$ws.Range("L132").Value = 14494.9995
$ws.Range("M132").Value = -11329.1
$ws.Range("N132").Value = -19554.9995

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3588.4546
$ws.Range("I136").Value = 3588.4546
$ws.Range("K136").Value = 10765.3638
$ws.Range("M136").Value = -8215.363799999999

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4474.2856
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 4830
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 4830
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -6078

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4474.2856
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 4830
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 24150
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -30390

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3423
$ws.Range("I132").Value = 2744.375
$ws.Range("J132").Value = 6137.5
$ws.Range("K132").Value = 8233.125
$ws.Range("L132").Value = 18412.5
$ws.Range("M132").Value = -5703.125
$ws.Range("N132").Value = -23472.5
